$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "1. Open https://abantecart.codifyme.co.nz`n2. Check if logo element is displayed`n3. Check if logo image exists in the specified source Url"
$ws.Range("F5").Value = "2. Logo element should be displayed`n3. Logo image file exists in the specified source Url."
